$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "U" column of values in D2:D6
$ws.Range("D2").Value = "U"
$ws.Range("D3").Value = "U"
$ws.Range("D4").Value = "U"
$ws.Range("D5").Value = "U"
$ws.Range("D6").Value = "U"

# Remove the old F and G columns' data (rows 2-6)
$ws.Range("F2:G6").ClearContents()

# Update the active selection to D7
$ws.Range("D7").Select()
